$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 980636.2
$ws.Range("I9").Value = 1225677.2
$ws.Range("K9").Value = 1225677.2
$ws.Range("M9").Value = -1225508.2

$ws.Range("H32").Value = 13892055
$ws.Range("J32").Value = 4332.3335
$ws.Range("L32").Value = 4332.3335
$ws.Range("N32").Value = -4984.3335

$ws.Range("H53").Value = 1269.9286
$ws.Range("I53").Value = 1522.1
$ws.Range("J53").Value = 639.5
$ws.Range("K53").Value = 1522.1
$ws.Range("L53").Value = 639.5
$ws.Range("M53").Value = -885.0999999999999
$ws.Range("N53").Value = -1913.5

$ws.Range("H62").Value = 5995.5
$ws.Range("J62").Value = 5995
$ws.Range("L62").Value = 5995
$ws.Range("N62").Value = -7243

$ws.Range("H65").Value = 5995.5
$ws.Range("J65").Value = 5995
$ws.Range("L65").Value = 29975
$ws.Range("N65").Value = -36215

$ws.Range("H86").Value = 3654.111
$ws.Range("I86").Value = 2860.875
$ws.Range("K86").Value = 2860.875
$ws.Range("M86").Value = -1737.875

$ws.Range("H89").Value = 3654.111
$ws.Range("I89").Value = 2860.875
$ws.Range("K89").Value = 14304.375
$ws.Range("M89").Value = -8688.375

$ws.Range("H97").Value = 3677.5
$ws.Range("J97").Value = 3677.5
$ws.Range("L97").Value = 11032.5
$ws.Range("N97").Value = -12024.5

$ws.Range("H116").Value = 10706.857
$ws.Range("I116").Value = 2649.3333
$ws.Range("K116").Value = 2649.3333
$ws.Range("M116").Value = 792.6667000000002

$ws.Range("H132").Value = 4368.8237
$ws.Range("I132").Value = 4657.839
$ws.Range("K132").Value = 13973.517
$ws.Range("M132").Value = -11443.517

$ws.Range("H135").Value = 933
$ws.Range("I135").Value = 466.1
$ws.Range("K135").Value = 4194.900000000001
$ws.Range("M135").Value = -1659.900000000001

$ws.Range("H137").Value = 6477.2
$ws.Range("I137").Value = 2734.913
$ws.Range("J137").Value = 18773.285
$ws.Range("K137").Value = 8204.739
$ws.Range("L137").Value = 56319.855
$ws.Range("M137").Value = -5654.739
$ws.Range("N137").Value = -61419.855

$ws.Range("H141").Value = 2545.3157
$ws.Range("I141").Value = 2403.111
$ws.Range("K141").Value = 7209.333
$ws.Range("M141").Value = -2029.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1896.7805
$ws.Range("I32").Value = 1919.2
$ws.Range("K32").Value = 1919.2
$ws.Range("M32").Value = -1632.2

$ws.Range("H61").Value = 2162
$ws.Range("J61").Value = 4434.8887
$ws.Range("L61").Value = 4434.8887
$ws.Range("N61").Value = -4858.8887

$ws.Range("H74").Value = 165238.8
$ws.Range("I74").Value = 232880
$ws.Range("K74").Value = 232880
$ws.Range("M74").Value = -232006

$ws.Range("H77").Value = 165238.8
$ws.Range("I77").Value = 232880
$ws.Range("K77").Value = 1164400
$ws.Range("M77").Value = -1160032

$ws.Range("H122").Value = 2743.3547
$ws.Range("I122").Value = 2834.4814
$ws.Range("K122").Value = 8503.4442
$ws.Range("M122").Value = -6053.4442

$ws.Range("H132").Value = 2194.5952
$ws.Range("I132").Value = 1975.9714
$ws.Range("J132").Value = 3287.7144
$ws.Range("K132").Value = 5927.914199999999
$ws.Range("L132").Value = 9863.143199999999
$ws.Range("M132").Value = -3397.914199999999
$ws.Range("N132").Value = -14923.1432

$ws.Range("H136").Value = 2162
$ws.Range("J136").Value = 4434.8887
$ws.Range("L136").Value = 13304.6661
$ws.Range("N136").Value = -18404.6661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2791.923
$ws.Range("I99").Value = 2225
$ws.Range("K99").Value = 2225
$ws.Range("M99").Value = -727

$ws.Range("H134").Value = 1963.7391
$ws.Range("I134").Value = 1564.7778
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 4694.3334
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -2159.3334
$ws.Range("N134").Value = -15270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 90136.37
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 90136.37
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 90136.37
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -90360.37

$ws.Range("H13").Value = 5399
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40490

$ws.Range("H31").Value = 3293228.2
$ws.Range("I31").Value = 2702.7097
$ws.Range("J31").Value = 17865556
$ws.Range("K31").Value = 2702.7097
$ws.Range("L31").Value = 17865556
$ws.Range("M31").Value = -2407.7097
$ws.Range("N31").Value = -17866146

$ws.Range("H34").Value = 3293228.2
$ws.Range("I34").Value = 2702.7097
$ws.Range("J34").Value = 17865556
$ws.Range("K34").Value = 2702.7097
$ws.Range("L34").Value = 17865556
$ws.Range("M34").Value = -2500.7097
$ws.Range("N34").Value = -17865960

$ws.Range("H122").Value = 391.95
$ws.Range("I122").Value = 348.4
$ws.Range("K122").Value = 1045.2
$ws.Range("M122").Value = 1404.8

$ws.Range("H132").Value = 3993.658
$ws.Range("I132").Value = 3674.7307
$ws.Range("J132").Value = 4684.6665
$ws.Range("K132").Value = 11024.1921
$ws.Range("L132").Value = 14053.9995
$ws.Range("M132").Value = -8494.1921
$ws.Range("N132").Value = -19113.9995

$ws.Range("H134").Value = 4189.25
$ws.Range("I134").Value = 4248.7915
$ws.Range("K134").Value = 12746.3745
$ws.Range("M134").Value = -10211.3745

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2890.3635
$ws.Range("I3").Value = 2099.25
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 6297.75
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -6185.75
$ws.Range("N3").Value = -15224

$ws.Range("H122").Value = 876.5
$ws.Range("I122").Value = 502
$ws.Range("K122").Value = 4518
$ws.Range("M122").Value = -2068

$ws.Range("H131").Value = 5893.7144
$ws.Range("J131").Value = 2452.875
$ws.Range("L131").Value = 7358.625
$ws.Range("N131").Value = -17438.625

$ws.Range("H132").Value = 2171
$ws.Range("J132").Value = 2467
$ws.Range("L132").Value = 22203
$ws.Range("N132").Value = -27263

$ws.Range("H133").Value = 3000

$ws.Range("H140").Value = 3456.1
$ws.Range("I140").Value = 3456.1
$ws.Range("K140").Value = 10368.3
$ws.Range("M140").Value = -5188.299999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 250002900
$ws.Range("I80").Value = 333335040
$ws.Range("J80").Value = 6500
$ws.Range("K80").Value = 333335040
$ws.Range("L80").Value = 6500
$ws.Range("M80").Value = -333334042
$ws.Range("N80").Value = -8496

$ws.Range("H83").Value = 250002900
$ws.Range("I83").Value = 333335040
$ws.Range("J83").Value = 6500
$ws.Range("K83").Value = 1666675200
$ws.Range("L83").Value = 32500
$ws.Range("M83").Value = -1666670208
$ws.Range("N83").Value = -42484

$ws.Range("H97").Value = 1693.9
$ws.Range("I97").Value = 1682.1111
$ws.Range("J97").Value = 1800
$ws.Range("K97").Value = 1682.1111
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -1186.1111
$ws.Range("N97").Value = -2792

$ws.Range("H102").Value = 1683.9412
$ws.Range("I102").Value = 1173.6154
$ws.Range("J102").Value = 1999.8572
$ws.Range("K102").Value = 1173.6154
$ws.Range("L102").Value = 1999.8572
$ws.Range("M102").Value = 448.3846000000001
$ws.Range("N102").Value = -5243.8572

$ws.Range("H132").Value = 1953.9546
$ws.Range("I132").Value = 1489
$ws.Range("J132").Value = 2565.7368
$ws.Range("K132").Value = 4467
$ws.Range("L132").Value = 7697.2104
$ws.Range("M132").Value = -1937
$ws.Range("N132").Value = -12757.2104

$ws.Range("H141").Value = 21000.5
$ws.Range("J141").Value = 21000.5
$ws.Range("L141").Value = 21000.5
$ws.Range("N141").Value = -31360.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19708.334
$ws.Range("I40").Value = 19708.334
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 19708.334
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -19572.334
$ws.Range("N40").ClearContents()

$ws.Range("H132").Value = 4373
$ws.Range("I132").Value = 3656.1667
$ws.Range("K132").Value = 10968.5001
$ws.Range("M132").Value = -8438.500100000001

$ws.Range("H136").Value = 5116
$ws.Range("I136").Value = 2248.6667
$ws.Range("J136").Value = 7983.3335
$ws.Range("K136").Value = 6746.000100000001
$ws.Range("L136").Value = 23950.0005
$ws.Range("M136").Value = -4196.000100000001
$ws.Range("N136").Value = -29050.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H132").Value = 1272.8379
$ws.Range("I132").Value = 1173.6207
$ws.Range("K132").Value = 3520.8621
$ws.Range("M132").Value = -990.8620999999998

$ws.Range("H136").Value = 6368.5
$ws.Range("I136").Value = 3459.4
$ws.Range("K136").Value = 10378.2
$ws.Range("M136").Value = -7828.200000000001
